# Generate Report for Handoff
# Adds two new localization entries (8fe14779-... and d8af8ba1-...) to the
# Overview / zh-cn / de-de worksheets, sorted alphabetically by file name,
# and refreshes the handoff timestamps.

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/3d697fc6072c8373828f7c1ad0413b3142307908/e2e"
$zhBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f519891cbc541a034c2b43d35dca49aca8d8d6e8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/high"
$deBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/34240f6c40914bf0b2b037e6b00ab28cad582347/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/high"

$file0 = "0c12b682-acc1-4ea5-aedd-a14dcb51d269"
$file1 = "8fe14779-40ee-46c9-906b-a989f70d66c9"
$file2 = "cf7994e8-3ca9-439c-b2ec-a01ab685158e"
$file3 = "d8af8ba1-0e07-4e29-b0d9-fb235654b05f"

$zhHash0 = "5d632f83fa835505597aaa80060172288f39d6ea"
$zhHash1 = "2bac596d149c23ddb1e96dc6949a0c97f5adddbc"
$zhHash2 = "a375b1bca50934f18cace6f778ecedb9d8ecb18b"
$zhHash3 = "33dc23fc83e6624b1d34f36750abfc9aaf8283bc"

$deHash0 = $zhHash0
$deHash1 = $zhHash1
$deHash2 = $zhHash2
$deHash3 = $zhHash3

$newHandoffDate = "2016-14-14 00:14:55"
$zhDatetime = "2016-03-14 00:14:51"
$deDatetime = "2016-03-14 00:14:55"
$handbackDatetime = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()

$ws1.Range("D2").Value2 = $newHandoffDate
$ws1.Range("D3").Value2 = $newHandoffDate

$ws1.Range("A4").Value2 = $file2 + ".md"
$ws1.Range("B4").Value2 = "Ready for handoff"
$ws1.Range("C4").Value2 = "Ready for handoff"
$ws1.Range("D4").Value2 = $newHandoffDate

$ws1.Range("A5").Value2 = $file3 + ".md"
$ws1.Range("B5").Value2 = "Ready for handoff"
$ws1.Range("C5").Value2 = "Ready for handoff"
$ws1.Range("D5").Value2 = $newHandoffDate

$ws1.Range("A3").Value2 = $file1 + ".md"
$ws1.Range("B3").Value2 = "Ready for handoff"
$ws1.Range("C3").Value2 = "Ready for handoff"
$ws1.Range("D3").Value2 = $newHandoffDate

$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdBase + "/" + $file0 + ".md", "", "", $file0 + ".md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdBase + "/" + $file1 + ".md", "", "", $file1 + ".md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), $mdBase + "/" + $file2 + ".md", "", "", $file2 + ".md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), $mdBase + "/" + $file3 + ".md", "", "", $file3 + ".md")

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()

$ws2.Range("E2").Value2 = $zhDatetime
$ws2.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("A4").Value2 = $file2 + ".md"
$ws2.Range("B4").Value2 = ".md"
$ws2.Range("C4").Value2 = "Ready for handoff"
$ws2.Range("D4").Value2 = $file2 + "." + $zhHash2 + ".zh-cn.xlf"
$ws2.Range("E4").Value2 = $zhDatetime
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H4").Value2 = $handbackDatetime
$ws2.Range("I4").Value2 = "Include"

$ws2.Range("A5").Value2 = $file3 + ".md"
$ws2.Range("B5").Value2 = ".md"
$ws2.Range("C5").Value2 = "Ready for handoff"
$ws2.Range("D5").Value2 = $file3 + "." + $zhHash3 + ".zh-cn.xlf"
$ws2.Range("E5").Value2 = $zhDatetime
$ws2.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H5").Value2 = $handbackDatetime
$ws2.Range("I5").Value2 = "Include"

$ws2.Range("A3").Value2 = $file1 + ".md"
$ws2.Range("B3").Value2 = ".md"
$ws2.Range("C3").Value2 = "Ready for handoff"
$ws2.Range("D3").Value2 = $file1 + "." + $zhHash1 + ".zh-cn.xlf"
$ws2.Range("E3").Value2 = $zhDatetime
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value2 = $handbackDatetime
$ws2.Range("I3").Value2 = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdBase + "/" + $file0 + ".md", "", "", $file0 + ".md")
$ws2.Hyperlinks.Add($ws2.Range("B2"), $mdBase + "/" + $file0 + ".md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), $zhBase + "/" + $file0 + "." + $zhHash0 + ".zh-cn.xlf", "", "", $file0 + "." + $zhHash0 + ".zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdBase + "/" + $file1 + ".md", "", "", $file1 + ".md")
$ws2.Hyperlinks.Add($ws2.Range("B3"), $mdBase + "/" + $file1 + ".md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), $zhBase + "/" + $file1 + "." + $zhHash1 + ".zh-cn.xlf", "", "", $file1 + "." + $zhHash1 + ".zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), $mdBase + "/" + $file2 + ".md", "", "", $file2 + ".md")
$ws2.Hyperlinks.Add($ws2.Range("B4"), $mdBase + "/" + $file2 + ".md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), $zhBase + "/" + $file2 + "." + $zhHash2 + ".zh-cn.xlf", "", "", $file2 + "." + $zhHash2 + ".zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A5"), $mdBase + "/" + $file3 + ".md", "", "", $file3 + ".md")
$ws2.Hyperlinks.Add($ws2.Range("B5"), $mdBase + "/" + $file3 + ".md", "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), $zhBase + "/" + $file3 + "." + $zhHash3 + ".zh-cn.xlf", "", "", $file3 + "." + $zhHash3 + ".zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()

$ws3.Range("E2").Value2 = $deDatetime
$ws3.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Range("A4").Value2 = $file2 + ".md"
$ws3.Range("B4").Value2 = ".md"
$ws3.Range("C4").Value2 = "Ready for handoff"
$ws3.Range("D4").Value2 = $file2 + "." + $deHash2 + ".de-de.xlf"
$ws3.Range("E4").Value2 = $deDatetime
$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H4").Value2 = $handbackDatetime
$ws3.Range("I4").Value2 = "Include"

$ws3.Range("A5").Value2 = $file3 + ".md"
$ws3.Range("B5").Value2 = ".md"
$ws3.Range("C5").Value2 = "Ready for handoff"
$ws3.Range("D5").Value2 = $file3 + "." + $deHash3 + ".de-de.xlf"
$ws3.Range("E5").Value2 = $deDatetime
$ws3.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H5").Value2 = $handbackDatetime
$ws3.Range("I5").Value2 = "Include"

$ws3.Range("A3").Value2 = $file1 + ".md"
$ws3.Range("B3").Value2 = ".md"
$ws3.Range("C3").Value2 = "Ready for handoff"
$ws3.Range("D3").Value2 = $file1 + "." + $deHash1 + ".de-de.xlf"
$ws3.Range("E3").Value2 = $deDatetime
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").Value2 = $handbackDatetime
$ws3.Range("I3").Value2 = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdBase + "/" + $file0 + ".md", "", "", $file0 + ".md")
$ws3.Hyperlinks.Add($ws3.Range("B2"), $mdBase + "/" + $file0 + ".md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), $deBase + "/" + $file0 + "." + $deHash0 + ".de-de.xlf", "", "", $file0 + "." + $deHash0 + ".de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdBase + "/" + $file1 + ".md", "", "", $file1 + ".md")
$ws3.Hyperlinks.Add($ws3.Range("B3"), $mdBase + "/" + $file1 + ".md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), $deBase + "/" + $file1 + "." + $deHash1 + ".de-de.xlf", "", "", $file1 + "." + $deHash1 + ".de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), $mdBase + "/" + $file2 + ".md", "", "", $file2 + ".md")
$ws3.Hyperlinks.Add($ws3.Range("B4"), $mdBase + "/" + $file2 + ".md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), $deBase + "/" + $file2 + "." + $deHash2 + ".de-de.xlf", "", "", $file2 + "." + $deHash2 + ".de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A5"), $mdBase + "/" + $file3 + ".md", "", "", $file3 + ".md")
$ws3.Hyperlinks.Add($ws3.Range("B5"), $mdBase + "/" + $file3 + ".md", "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), $deBase + "/" + $file3 + "." + $deHash3 + ".de-de.xlf", "", "", $file3 + "." + $deHash3 + ".de-de.xlf")

Write-Host "Report regenerated for handoff."
